# Weekly fruit/vegetable update: a new weekly record was inserted at row 334,
# pushing the existing rows 334-353 down to 335-354.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 334 (shifts rows 334:353 down to 335:354)
$ws.Rows.Item(334).Insert()

# Populate the new row 334 with the new weekly record
$ws.Cells.Item(334, 1).Value  = 9
$ws.Cells.Item(334, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(334, 3).Value  = "Metropolitana"
$ws.Cells.Item(334, 4).Value  = 44610
$ws.Cells.Item(334, 5).Value  = 13
$ws.Cells.Item(334, 6).Value  = 100112039
$ws.Cells.Item(334, 7).Value  = "Ciboulette"
$ws.Cells.Item(334, 8).Value  = "Sin especificar"
$ws.Cells.Item(334, 9).Value  = "Primera"
$ws.Cells.Item(334, 10).Value = 250
$ws.Cells.Item(334, 11).Value = 1000
$ws.Cells.Item(334, 12).Value = 1200
$ws.Cells.Item(334, 13).Value = 1100
$ws.Cells.Item(334, 14).Value = "`$/docena de atados"
$ws.Cells.Item(334, 15).Value = "Región Metropolitana"
$ws.Cells.Item(334, 16).Value = 367
$ws.Cells.Item(334, 17).Value = 3
$ws.Cells.Item(334, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by the
# rest of column D.
$ws.Cells.Item(334, 4).NumberFormat = $ws.Cells.Item(335, 4).NumberFormat()
